$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(13, 8).Value = 1504  # H13: was 66666830
$ws.Cells.Item(13, 9).Value = 0  # I13: was 100000000
$ws.Cells.Item(13, 10).Value = 1504  # J13: was 500
$ws.Cells.Item(13, 11).Value = 0  # K13: was 100000000
$ws.Cells.Item(13, 12).Value = 1504  # L13: was 500
$ws.Cells.Item(13, 13).ClearContents()  # M13: was -99999856
$ws.Cells.Item(13, 14).Value = -1792  # N13: was -788

$ws.Cells.Item(26, 8).Value = 4559.857  # H26: was 4500
$ws.Cells.Item(26, 9).Value = 730  # I26: was 3375
$ws.Cells.Item(26, 10).Value = 9666.333  # J26: was 5062.5
$ws.Cells.Item(26, 11).Value = 730  # K26: was 3375
$ws.Cells.Item(26, 12).Value = 9666.333  # L26: was 5062.5
$ws.Cells.Item(26, 13).Value = -400  # M26: was -3045
$ws.Cells.Item(26, 14).Value = -10326.333  # N26: was -5722.5

$ws.Cells.Item(28, 8).Value = 6919.2354  # H28: was 10290.533
$ws.Cells.Item(28, 9).Value = 4555.231  # I28: was 7121.727
$ws.Cells.Item(28, 10).Value = 14602.25  # J28: was 19004.75
$ws.Cells.Item(28, 11).Value = 4555.231  # K28: was 7121.727
$ws.Cells.Item(28, 12).Value = 14602.25  # L28: was 19004.75
$ws.Cells.Item(28, 13).Value = -4363.231  # M28: was -6929.727
$ws.Cells.Item(28, 14).Value = -14986.25  # N28: was -19388.75

$ws.Cells.Item(61, 8).Value = 591066.2  # H61: was 591066.5
$ws.Cells.Item(61, 9).Value = 502121.7  # I61: was 528501.8
$ws.Cells.Item(61, 10).Value = 718129.8  # J61: was 670315.2
$ws.Cells.Item(61, 11).Value = 502121.7  # K61: was 528501.8
$ws.Cells.Item(61, 12).Value = 718129.8  # L61: was 670315.2
$ws.Cells.Item(61, 13).Value = -501909.7  # M61: was -528289.8
$ws.Cells.Item(61, 14).Value = -718553.8  # N61: was -670739.2

$ws.Cells.Item(74, 8).Value = 12322505  # H74: was 11685147
$ws.Cells.Item(74, 9).Value = 9035746  # I74: was 8374612
$ws.Cells.Item(74, 11).Value = 9035746  # K74: was 8374612
$ws.Cells.Item(74, 13).Value = -9034872  # M74: was -8373738

$ws.Cells.Item(77, 8).Value = 12322505  # H77: was 11685147
$ws.Cells.Item(77, 9).Value = 9035746  # I77: was 8374612
$ws.Cells.Item(77, 11).Value = 45178730  # K77: was 41873060
$ws.Cells.Item(77, 13).Value = -45174362  # M77: was -41868692

$ws.Cells.Item(93, 8).Value = 25000  # H93: was 26000
$ws.Cells.Item(93, 10).Value = 25000  # J93: was 26000
$ws.Cells.Item(93, 12).Value = 25000  # L93: was 26000
$ws.Cells.Item(93, 14).Value = -29992  # N93: was -30992

$ws.Cells.Item(99, 8).Value = 6919.2354  # H99: was 10290.533
$ws.Cells.Item(99, 9).Value = 4555.231  # I99: was 7121.727
$ws.Cells.Item(99, 10).Value = 14602.25  # J99: was 19004.75
$ws.Cells.Item(99, 11).Value = 4555.231  # K99: was 7121.727
$ws.Cells.Item(99, 12).Value = 14602.25  # L99: was 19004.75
$ws.Cells.Item(99, 13).Value = -1560.231  # M99: was -4126.727
$ws.Cells.Item(99, 14).Value = -20592.25  # N99: was -24994.75

$ws.Cells.Item(132, 8).Value = 1987.5818  # H132: was 2109.4468
$ws.Cells.Item(132, 9).Value = 1045.3226  # I132: was 1182.5385
$ws.Cells.Item(132, 10).Value = 3204.6667  # J132: was 3257.0476
$ws.Cells.Item(132, 11).Value = 3135.9678  # K132: was 3547.6155
$ws.Cells.Item(132, 12).Value = 9614.000100000001  # L132: was 9771.1428
$ws.Cells.Item(132, 13).Value = -605.9677999999999  # M132: was -1017.6155
$ws.Cells.Item(132, 14).Value = -14674.0001  # N132: was -14831.1428

$ws.Cells.Item(136, 8).Value = 591066.2  # H136: was 591066.5
$ws.Cells.Item(136, 9).Value = 502121.7  # I136: was 528501.8
$ws.Cells.Item(136, 10).Value = 718129.8  # J136: was 670315.2
$ws.Cells.Item(136, 11).Value = 1506365.1  # K136: was 1585505.4
$ws.Cells.Item(136, 12).Value = 2154389.4  # L136: was 2010945.6
$ws.Cells.Item(136, 13).Value = -1503815.1  # M136: was -1582955.4
$ws.Cells.Item(136, 14).Value = -2159489.4  # N136: was -2016045.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(21, 8).Value = 22689  # H21: was 23000
$ws.Cells.Item(21, 10).Value = 22689  # J21: was 23000
$ws.Cells.Item(21, 12).Value = 22689  # L21: was 23000
$ws.Cells.Item(21, 14).Value = -23161  # N21: was -23472

$ws.Cells.Item(26, 8).Value = 22000  # H26: was 24666.666
$ws.Cells.Item(26, 10).Value = 29000  # J26: was 29500
$ws.Cells.Item(26, 12).Value = 29000  # L26: was 29500
$ws.Cells.Item(26, 14).Value = -29584  # N26: was -30084

$ws.Cells.Item(40, 8).Value = 26000  # H40: was 25974
$ws.Cells.Item(40, 10).Value = 26000  # J40: was 25974
$ws.Cells.Item(40, 12).Value = 26000  # L40: was 25974
$ws.Cells.Item(40, 14).Value = -26530  # N40: was -26504

$ws.Cells.Item(96, 8).Value = 16420  # H96: was 11378
$ws.Cells.Item(96, 9).Value = 3840  # I96: was 2317
$ws.Cells.Item(96, 10).Value = 29000  # J96: was 29500
$ws.Cells.Item(96, 11).Value = 3840  # K96: was 2317
$ws.Cells.Item(96, 12).Value = 29000  # L96: was 29500
$ws.Cells.Item(96, 13).Value = -1094  # M96: was 429
$ws.Cells.Item(96, 14).Value = -34492  # N96: was -34992

$ws.Cells.Item(98, 8).Value = 32333  # H98: was 29000
$ws.Cells.Item(98, 10).Value = 32333  # J98: was 29000
$ws.Cells.Item(98, 12).Value = 32333  # L98: was 29000
$ws.Cells.Item(98, 14).Value = -38323  # N98: was -34990

$ws.Cells.Item(134, 8).Value = 6424.433  # H134: was 5748.8
$ws.Cells.Item(134, 9).Value = 7197.087  # I134: was 6896.1665
$ws.Cells.Item(134, 10).Value = 3885.7144  # J134: was 3245.4546
$ws.Cells.Item(134, 11).Value = 21591.261  # K134: was 20688.4995
$ws.Cells.Item(134, 12).Value = 11657.1432  # L134: was 9736.3638
$ws.Cells.Item(134, 13).Value = -19056.261  # M134: was -18153.4995
$ws.Cells.Item(134, 14).Value = -16727.1432  # N134: was -14806.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 7795.8335  # H58: was 8582.875
$ws.Cells.Item(58, 9).Value = 11401.1  # I58: was 12556.889
$ws.Cells.Item(58, 10).Value = 3289.25  # J58: was 3473.4285
$ws.Cells.Item(58, 11).Value = 11401.1  # K58: was 12556.889
$ws.Cells.Item(58, 12).Value = 3289.25  # L58: was 3473.4285
$ws.Cells.Item(58, 13).Value = -11198.1  # M58: was -12353.889
$ws.Cells.Item(58, 14).Value = -3695.25  # N58: was -3879.4285

$ws.Cells.Item(99, 8).Value = 85133.5  # H99: was 49575.668
$ws.Cells.Item(99, 9).Value = 101560.2  # I99: was 64351.188
$ws.Cells.Item(99, 10).Value = 3000  # J99: was 2294
$ws.Cells.Item(99, 11).Value = 101560.2  # K99: was 64351.188
$ws.Cells.Item(99, 12).Value = 3000  # L99: was 2294
$ws.Cells.Item(99, 13).Value = -100062.2  # M99: was -62853.188
$ws.Cells.Item(99, 14).Value = -5996  # N99: was -5290

$ws.Cells.Item(114, 8).Value = 16660  # H114: was 19666.334
$ws.Cells.Item(114, 10).Value = 16660  # J114: was 19666.334
$ws.Cells.Item(114, 12).Value = 16660  # L114: was 19666.334
$ws.Cells.Item(114, 14).Value = -25338  # N114: was -28344.334

$ws.Cells.Item(126, 8).Value = 85133.5  # H126: was 49575.668
$ws.Cells.Item(126, 9).Value = 101560.2  # I126: was 64351.188
$ws.Cells.Item(126, 10).Value = 3000  # J126: was 2294
$ws.Cells.Item(126, 11).Value = 304680.6  # K126: was 193053.564
$ws.Cells.Item(126, 12).Value = 9000  # L126: was 6882
$ws.Cells.Item(126, 13).Value = -302210.6  # M126: was -190583.564
$ws.Cells.Item(126, 14).Value = -13940  # N126: was -11822

$ws.Cells.Item(132, 8).Value = 14708218  # H132: was 13891136
$ws.Cells.Item(132, 9).Value = 25001602  # I132: was 22728796
$ws.Cells.Item(132, 11).Value = 75004806  # K132: was 68186388
$ws.Cells.Item(132, 13).Value = -75002276  # M132: was -68183858

$ws.Cells.Item(134, 8).Value = 16668140  # H134: was 31252004
$ws.Cells.Item(134, 9).Value = 25000806  # I134: was 45455490
$ws.Cells.Item(134, 10).Value = 2809.1  # J134: was 4342.8
$ws.Cells.Item(134, 11).Value = 75002418  # K134: was 136366470
$ws.Cells.Item(134, 12).Value = 8427.3  # L134: was 13028.4
$ws.Cells.Item(134, 13).Value = -74999883  # M134: was -136363935
$ws.Cells.Item(134, 14).Value = -13497.3  # N134: was -18098.4

$ws.Cells.Item(136, 8).Value = 7795.8335  # H136: was 8582.875
$ws.Cells.Item(136, 9).Value = 11401.1  # I136: was 12556.889
$ws.Cells.Item(136, 10).Value = 3289.25  # J136: was 3473.4285
$ws.Cells.Item(136, 11).Value = 34203.3  # K136: was 37670.667
$ws.Cells.Item(136, 12).Value = 9867.75  # L136: was 10420.2855
$ws.Cells.Item(136, 13).Value = -31653.3  # M136: was -35120.667
$ws.Cells.Item(136, 14).Value = -14967.75  # N136: was -15520.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 2147.647  # H81: was 2266.3333
$ws.Cells.Item(81, 10).Value = 2347.6924  # J81: was 2545.9092
$ws.Cells.Item(81, 12).Value = 7043.0772  # L81: was 7637.7276
$ws.Cells.Item(81, 14).Value = -9289.0772  # N81: was -9883.7276

$ws.Cells.Item(82, 8).Value = 7771.4287  # H82: was 7000
$ws.Cells.Item(82, 10).Value = 8900  # J82: was 9000
$ws.Cells.Item(82, 12).Value = 26700  # L82: was 27000
$ws.Cells.Item(82, 14).Value = -27512  # N82: was -27812

$ws.Cells.Item(84, 8).Value = 2147.647  # H84: was 2266.3333
$ws.Cells.Item(84, 10).Value = 2347.6924  # J84: was 2545.9092
$ws.Cells.Item(84, 12).Value = 21129.2316  # L84: was 22913.1828
$ws.Cells.Item(84, 14).Value = -32361.2316  # N84: was -34145.1828

$ws.Cells.Item(85, 8).Value = 7771.4287  # H85: was 7000
$ws.Cells.Item(85, 10).Value = 8900  # J85: was 9000
$ws.Cells.Item(85, 12).Value = 26700  # L85: was 27000
$ws.Cells.Item(85, 14).Value = -29508  # N85: was -29808

$ws.Cells.Item(131, 8).Value = 2450  # H131: was 1350.641
$ws.Cells.Item(131, 9).Value = 1272.2222  # I131: was 1274.4445
$ws.Cells.Item(131, 10).Value = 7750  # J131: was 1373.5
$ws.Cells.Item(131, 11).Value = 3816.6666  # K131: was 3823.3335
$ws.Cells.Item(131, 12).Value = 23250  # L131: was 4120.5
$ws.Cells.Item(131, 13).Value = 1223.3334  # M131: was 1216.6665
$ws.Cells.Item(131, 14).Value = -33330  # N131: was -14200.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(94, 8).Value = 22000  # H94: was 22333.334
$ws.Cells.Item(94, 10).Value = 22000  # J94: was 22333.334
$ws.Cells.Item(94, 12).Value = 22000  # L94: was 22333.334
$ws.Cells.Item(94, 14).Value = -23352  # N94: was -23685.334

$ws.Cells.Item(100, 8).Value = 37339  # H100: was 0
$ws.Cells.Item(100, 10).Value = 37339  # J100: was 0
$ws.Cells.Item(100, 12).Value = 37339  # L100: was 0
$ws.Cells.Item(100, 14).Value = -39503  # N100: was None

$ws.Cells.Item(122, 8).Value = 774.73334  # H122: was 1483.3334
$ws.Cells.Item(122, 9).Value = 800.7857  # I122: was 1483.3334
$ws.Cells.Item(122, 10).Value = 410  # J122: was 0
$ws.Cells.Item(122, 11).Value = 2402.3571  # K122: was 4450.0002
$ws.Cells.Item(122, 12).Value = 1230  # L122: was 0
$ws.Cells.Item(122, 13).Value = 47.64289999999983  # M122: was -2000.0002
$ws.Cells.Item(122, 14).Value = -6130  # N122: was None

$ws.Cells.Item(132, 8).Value = 3128722.8  # H132: was 4036390
$ws.Cells.Item(132, 9).Value = 4811058  # I132: was 5685645.5
$ws.Cells.Item(132, 10).Value = 4385.7144  # J132: was 4876.222
$ws.Cells.Item(132, 11).Value = 14433174  # K132: was 17056936.5
$ws.Cells.Item(132, 12).Value = 13157.1432  # L132: was 14628.666
$ws.Cells.Item(132, 13).Value = -14430644  # M132: was -17054406.5
$ws.Cells.Item(132, 14).Value = -18217.1432  # N132: was -19688.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 3187.0159  # H136: was 3524.5
$ws.Cells.Item(136, 9).Value = 1660.921  # I136: was 1832.9412
$ws.Cells.Item(136, 10).Value = 5506.68  # J136: was 6138.727
$ws.Cells.Item(136, 11).Value = 4982.763  # K136: was 5498.8236
$ws.Cells.Item(136, 12).Value = 16520.04  # L136: was 18416.181
$ws.Cells.Item(136, 13).Value = -2432.763  # M136: was -2948.8236
$ws.Cells.Item(136, 14).Value = -21620.04  # N136: was -23516.181
